$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "242.83"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.99"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.411"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05959"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.427"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.501"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8140"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9201"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1436"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07387"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03302"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03080"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09350"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001570"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04707"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.005876"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "BitKan"
$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.001262"
$ws.Range("E19").Value = "18BitKanKAN"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004846"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.00007999"
$ws.Range("E21").Value = "20NitroExNTXWorstin24h"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.570"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.133"
$ws.Range("E23").Value = "22BTSETokenBTSE"
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0005889"
$ws.Range("E24").Value = "23OneONE"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3236"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1328"
$ws.Range("E27").Value = "26UpBotsUBXTBestin24h"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03948"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006339"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003899"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008912"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005180"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6999"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002140"
